$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$oldUrl = "https://github.com/vdelam00/SIBIVideojuegos"
$found = $tr.Find($oldUrl)

$start = $found.Start
$len = $found.Length

# Split the hyperlinked run "https://github.com/vdelam00/SIBIVideojuegos" into
# two runs - "https://" and "github.com/vdelam00/VideojuegosSIBI" - both keeping
# the same hyperlink (rId5) and formatting.
$schemePart = $tr.Characters($start, 8)
$domainPart = $tr.Characters($start + 8, $len - 8)
$domainPart.Text = "github.com/vdelam00/VideojuegosSIBI"

# The trailing single space run right after the link becomes two spaces.
$spaceRun = $tr.Characters($start + $len, 1)
$spaceRun.Text = "  "
